$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal:" value from "EP-3" to "EF-7,EP-3"
$ws.Range("B9").Value = "EF-7,EP-3"
$ws.Range("C9").Value = "EF-7,EP-3"

# Remove the "Requisitos:" rows (22 and 23)
$ws.Range("A22:C23").EntireRow.Delete()
